# Updated cryptos list (GitHub Actions price refresh).
# Price cells that look like a plain decimal number (e.g. "97.68") are
# written with a leading apostrophe so Excel stores them as text instead
# of auto-converting to a numeric value, matching the original inlineStr
# cell type. Prices with two dots (e.g. "43.770.75") or percent cells
# (e.g. "  +0.28%  ") are never auto-numeric, so they're set directly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.770.75"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "2.317.11"
$ws.Range("E3").Value = "  +4.10%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'97.68"
$ws.Range("E5").Value = "  +4.53%  "
$ws.Range("D6").Value = "'270.84"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "'0.626"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.625"
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("D10").Value = "'45.47"
$ws.Range("E10").Value = "  -0.43%  "
$ws.Range("D11").Value = "'0.0944"
$ws.Range("E11").Value = "  -2.97%  "
$ws.Range("D12").Value = "'8.03"
$ws.Range("E12").Value = "  -2.85%  "
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").Value = "2.656.10"
$ws.Range("E14").Value = "  +3.81%  "
$ws.Range("D15").Value = "'15.46"
$ws.Range("E15").Value = "  +3.40%  "
$ws.Range("D16").Value = "'0.869"
$ws.Range("E16").Value = "  +8.04%  "
$ws.Range("D17").Value = "2.315.16"
$ws.Range("E17").Value = "  +4.24%  "
$ws.Range("D18").Value = "43.725.09"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").Value = "'0.0000110"
$ws.Range("E19").Value = "  +3.54%  "
$ws.Range("D20").Value = "'6.39"
$ws.Range("E20").Value = "  +5.74%  "
$ws.Range("D21").Value = "'73.02"
$ws.Range("E21").Value = "  +3.49%  "
$ws.Range("D22").Value = "'239.34"
$ws.Range("E22").Value = "  +2.44%  "
$ws.Range("D23").Value = "'2.26"
$ws.Range("E23").Value = "  -4.02%  "
$ws.Range("D24").Value = "'9.41"
$ws.Range("E24").Value = "  +3.60%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").Value = "'2.53"
$ws.Range("E26").Value = "  +1.17%  "
$ws.Range("D27").Value = "'11.37"
$ws.Range("E27").Value = "  -0.46%  "
$ws.Range("E28").Value = "  -1.49%  "
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "'38.01"
$ws.Range("E30").Value = "  -10.50%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'22.37"
$ws.Range("E31").Value = "  +7.07%  "
$ws.Range("D32").Value = "'174.40"
$ws.Range("E32").Value = "  +1.15%  "
$ws.Range("D33").Value = "'0.0905"
$ws.Range("E33").Value = "  -1.49%  "
$ws.Range("D34").Value = "'5.47"
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("D35").Value = "'0.128"
$ws.Range("E35").Value = "  +2.88%  "
$ws.Range("D36").Value = "'0.0363"
$ws.Range("E36").Value = "  +3.29%  "
$ws.Range("D37").Value = "'0.109"
$ws.Range("E37").Value = "  -3.40%  "
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("E39").Value = "  -6.36%  "
$ws.Range("D40").Value = "'0.248"
$ws.Range("E40").Value = "  +12.22%  "
$ws.Range("E41").Value = "  +8.90%  "
$ws.Range("E42").Value = "  +18.02%  "
$ws.Range("D43").Value = "'12.16"
$ws.Range("E43").Value = "  -4.52%  "
$ws.Range("D44").Value = "'9.20"
$ws.Range("E44").Value = "  +9.80%  "
$ws.Range("D45").Value = "'62.10"
$ws.Range("E45").Value = "  -2.15%  "
$ws.Range("D46").Value = "'5.32"
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("E47").Value = "  +4.22%  "
$ws.Range("D48").Value = "'100.36"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("E50").Value = "  +15.38%  "
$ws.Range("D51").Value = "2.543.40"
$ws.Range("E51").Value = "  +3.97%  "

# Reset style on cells that were forced to text via apostrophe prefix,
# to avoid leaving a residual Text number-format style applied.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
